# This script updates the cryptocurrency price/volume snapshot cells
# (columns B-E, rows 2-51) on the active worksheet to match the latest
# scraped values, per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.063.61"
$ws.Range("E2").Value = "  +0.19%  "
# Row 3
$ws.Range("D3").Value = "2.959.00"
$ws.Range("E3").Value = "  +0.90%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
# Row 11
$ws.Range("E11").Value = "  -1.06%  "
# Row 12
$ws.Range("E12").Value = "  +2.31%  "
# Row 13
$ws.Range("D13").Value = "3.422.34"
$ws.Range("E13").Value = "  +0.82%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "
# Row 15
$ws.Range("E15").Value = "  +5.96%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +67.06%  "
# Row 17
$ws.Range("D17").Value = "2.963.12"
$ws.Range("E17").Value = "  +1.07%  "
# Row 18
$ws.Range("E18").Value = "  +2.91%  "
# Row 19
$ws.Range("D19").Value = "51.129.39"
$ws.Range("E19").Value = "  +0.35%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
# Row 22
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.86%  "
# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.59%  "
# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.80%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "267.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "
# Row 26
$ws.Range("E26").Value = "  -1.61%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.22%  "
# Row 28
$ws.Range("E28").Value = "  +0.00%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.167"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.12%  "
# Row 31
$ws.Range("E31").Value = "  -2.18%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.13%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.02%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "
# Row 35
$ws.Range("E35").Value = "  +2.30%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "
# Row 37
$ws.Range("E37").Value = "  +0.00%  "
# Row 38
$ws.Range("E38").Value = "  +10.52%  "
# Row 39
$ws.Range("E39").Value = "  +2.08%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.48%  "
# Row 41
$ws.Range("E41").Value = "  +1.77%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.41%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.76%  "
# Row 45
$ws.Range("E45").Value = "  +10.91%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.64%  "
# Row 47
$ws.Range("E47").Value = "  -0.73%  "
# Row 48
$ws.Range("D48").Value = "2.052.26"
$ws.Range("E48").Value = "  +4.24%  "
# Row 49
$ws.Range("E49").Value = "  -0.82%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.60%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.76%  "
